# Auto-generated script to apply scheduled-runner price updates to Sheets
# (Exodus_Profits workbook: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR tables)

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (27 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 99999
$ws.Range("J108").Value = 99999
$ws.Range("L108").Value = 99999
$ws.Range("N108").Value = -107679
$ws.Range("H109").Value = 99499.164
$ws.Range("J109").Value = 99499.164
$ws.Range("L109").Value = 99499.164
$ws.Range("N109").Value = -102273.164
$ws.Range("H123").Value = 77879.86
$ws.Range("J123").Value = 77879.86
$ws.Range("L123").Value = 77879.86
$ws.Range("N123").Value = -87679.86
$ws.Range("H134").Value = 54997.777
$ws.Range("J134").Value = 54997.777
$ws.Range("L134").Value = 54997.777
$ws.Range("N134").Value = -65137.777
$ws.Range("H137").Value = 364453.94
$ws.Range("I137").Value = 1714.4642
$ws.Range("K137").Value = 5143.392599999999
$ws.Range("M137").Value = -2593.392599999999
$ws.Range("H138").Value = 1855.45
$ws.Range("I138").Value = 1442.6666
$ws.Range("J138").Value = 2474.625
$ws.Range("K138").Value = 4327.9998
$ws.Range("L138").Value = 7423.875
$ws.Range("M138").Value = 812.0002000000004
$ws.Range("N138").Value = -17703.875

# --- Sheet: ARM (23 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6052.317
$ws.Range("I32").Value = 2742.422
$ws.Range("J32").Value = 15982
$ws.Range("K32").Value = 2742.422
$ws.Range("L32").Value = 15982
$ws.Range("M32").Value = -2455.422
$ws.Range("N32").Value = -16556
$ws.Range("H107").Value = 66482
$ws.Range("J107").Value = 66482
$ws.Range("L107").Value = 66482
$ws.Range("N107").Value = -74162
$ws.Range("H117").Value = 42346.668
$ws.Range("J117").Value = 42346.668
$ws.Range("L117").Value = 42346.668
$ws.Range("N117").Value = -51524.668
$ws.Range("H118").Value = 52997.145
$ws.Range("J118").Value = 52997.145
$ws.Range("L118").Value = 52997.145
$ws.Range("N118").Value = -56311.145
$ws.Range("H121").Value = 48649.066
$ws.Range("J121").Value = 48649.066
$ws.Range("L121").Value = 48649.066
$ws.Range("N121").Value = -52143.066

# --- Sheet: BSM (62 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4929.7915
$ws.Range("I86").Value = 3222.8333
$ws.Range("J86").Value = 10050.667
$ws.Range("K86").Value = 3222.8333
$ws.Range("L86").Value = 10050.667
$ws.Range("M86").Value = -2099.8333
$ws.Range("N86").Value = -12296.667
$ws.Range("H89").Value = 4929.7915
$ws.Range("I89").Value = 3222.8333
$ws.Range("J89").Value = 10050.667
$ws.Range("K89").Value = 16114.1665
$ws.Range("L89").Value = 50253.335
$ws.Range("M89").Value = -10498.1665
$ws.Range("N89").Value = -61485.335
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H107").Value = 1958.0952
$ws.Range("I107").Value = 1504.4667
$ws.Range("K107").Value = 1504.4667
$ws.Range("M107").Value = 415.5333000000001
$ws.Range("H108").Value = 94621.125
$ws.Range("J108").Value = 94621.125
$ws.Range("L108").Value = 94621.125
$ws.Range("N108").Value = -102301.125
$ws.Range("H109").Value = 74281.86
$ws.Range("J109").Value = 74281.86
$ws.Range("L109").Value = 74281.86
$ws.Range("N109").Value = -77055.86
$ws.Range("H110").Value = 83354.28999999999
$ws.Range("J110").Value = 83354.28999999999
$ws.Range("L110").Value = 83354.28999999999
$ws.Range("N110").Value = -91534.28999999999
$ws.Range("H114").Value = 90662.5
$ws.Range("J114").Value = 90662.5
$ws.Range("L114").Value = 90662.5
$ws.Range("N114").Value = -99340.5
$ws.Range("H118").Value = 72038
$ws.Range("J118").Value = 74786.664
$ws.Range("L118").Value = 74786.664
$ws.Range("N118").Value = -78100.664
$ws.Range("H122").Value = 72822.14
$ws.Range("J122").Value = 72822.14
$ws.Range("L122").Value = 72822.14
$ws.Range("N122").Value = -82622.14
$ws.Range("H132").Value = 28814.592
$ws.Range("J132").Value = 28814.592
$ws.Range("L132").Value = 28814.592
$ws.Range("N132").Value = -38934.592
$ws.Range("H134").Value = 4083.3103
$ws.Range("I134").Value = 3350.6667
$ws.Range("K134").Value = 10052.0001
$ws.Range("M134").Value = -7517.000100000001
$ws.Range("H135").Value = 52000
$ws.Range("J135").Value = 52000
$ws.Range("L135").Value = 52000
$ws.Range("N135").Value = -62140
$ws.Range("H138").Value = 99752.28999999999
$ws.Range("J138").Value = 99752.28999999999
$ws.Range("L138").Value = 99752.28999999999
$ws.Range("N138").Value = -110032.29

# --- Sheet: CRP (46 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3245.1667
$ws.Range("I31").Value = 2423.875
$ws.Range("K31").Value = 2423.875
$ws.Range("M31").Value = -2128.875
$ws.Range("H34").Value = 3245.1667
$ws.Range("I34").Value = 2423.875
$ws.Range("K34").Value = 2423.875
$ws.Range("M34").Value = -2221.875
$ws.Range("H108").Value = 53456.848
$ws.Range("I108").Value = 20000
$ws.Range("J108").Value = 56244.918
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 56244.918
$ws.Range("N108").Value = -63924.918
$ws.Range("M108").Value = -16160
$ws.Range("H114").Value = 39990.5
$ws.Range("J114").Value = 39990.5
$ws.Range("L114").Value = 39990.5
$ws.Range("N114").Value = -48668.5
$ws.Range("H116").Value = 83349.25
$ws.Range("J116").Value = 83349.25
$ws.Range("L116").Value = 83349.25
$ws.Range("N116").Value = -92527.25
$ws.Range("H117").Value = 38082
$ws.Range("J117").Value = 38082
$ws.Range("L117").Value = 38082
$ws.Range("N117").Value = -47260
$ws.Range("H118").Value = 99999
$ws.Range("J118").Value = 99999
$ws.Range("L118").Value = 99999
$ws.Range("N118").Value = -103313
$ws.Range("H119").Value = 94713.42999999999
$ws.Range("J119").Value = 94713.42999999999
$ws.Range("L119").Value = 94713.42999999999
$ws.Range("N119").Value = -104389.43
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258
$ws.Range("H134").Value = 25937.139
$ws.Range("I134").Value = 2895.1562
$ws.Range("J134").Value = 92968.37
$ws.Range("K134").Value = 8685.4686
$ws.Range("L134").Value = 278905.11
$ws.Range("M134").Value = -6150.4686
$ws.Range("N134").Value = -283975.11

# --- Sheet: CUL (18 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1863
$ws.Range("I36").Value = 1863
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5589
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -5420
$ws.Range("N36").Value = $null
$ws.Range("H107").Value = 823.8333
$ws.Range("I107").Value = 795.3333
$ws.Range("K107").Value = 2385.9999
$ws.Range("M107").Value = -465.9998999999998
$ws.Range("H129").Value = 1120.375
$ws.Range("I129").Value = 909.6667
$ws.Range("J129").Value = 1246.8
$ws.Range("K129").Value = 2729.0001
$ws.Range("L129").Value = 3740.4
$ws.Range("M129").Value = 2270.9999
$ws.Range("N129").Value = -13740.4

# --- Sheet: GSM (45 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 251.6
$ws.Range("I2").Value = 31.25
$ws.Range("J2").Value = 398.5
$ws.Range("K2").Value = 31.25
$ws.Range("L2").Value = 398.5
$ws.Range("M2").Value = 81.75
$ws.Range("N2").Value = -624.5
$ws.Range("H47").Value = 15304
$ws.Range("J47").Value = 15304
$ws.Range("L47").Value = 15304
$ws.Range("N47").Value = -16440
$ws.Range("H55").Value = 8008.75
$ws.Range("I55").Value = 7412
$ws.Range("J55").Value = 9003.333000000001
$ws.Range("K55").Value = 7412
$ws.Range("L55").Value = 9003.333000000001
$ws.Range("N55").Value = -9657.333000000001
$ws.Range("M55").Value = -7085
$ws.Range("H93").Value = 20463.334
$ws.Range("J93").Value = 20463.334
$ws.Range("L93").Value = 20463.334
$ws.Range("N93").Value = -24207.334
$ws.Range("H102").Value = 1495.8823
$ws.Range("I102").Value = 1327.3846
$ws.Range("J102").Value = 2043.5
$ws.Range("K102").Value = 1327.3846
$ws.Range("L102").Value = 2043.5
$ws.Range("M102").Value = 294.6153999999999
$ws.Range("N102").Value = -5287.5
$ws.Range("H110").Value = 67329.73
$ws.Range("J110").Value = 67329.73
$ws.Range("L110").Value = 67329.73
$ws.Range("N110").Value = -75509.73
$ws.Range("H114").Value = 60717.332
$ws.Range("J114").Value = 60717.332
$ws.Range("L114").Value = 60717.332
$ws.Range("N114").Value = -69395.33199999999
$ws.Range("H116").Value = 59996.668
$ws.Range("J116").Value = 59996.668
$ws.Range("L116").Value = 59996.668
$ws.Range("N116").Value = -69174.66800000001
$ws.Range("H132").Value = 6473.727
$ws.Range("I132").Value = 5618.5
$ws.Range("K132").Value = 16855.5
$ws.Range("M132").Value = -14325.5

# --- Sheet: LTW (8 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 7956.3125
$ws.Range("I55").Value = 1044.6666
$ws.Range("K55").Value = 1044.6666
$ws.Range("M55").Value = -871.6666
$ws.Range("H136").Value = 1905.6818
$ws.Range("I136").Value = 1594.5
$ws.Range("K136").Value = 4783.5
$ws.Range("M136").Value = -2233.5

# --- Sheet: WVR (15 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10393.625
$ws.Range("J41").Value = 10418.286
$ws.Range("L41").Value = 10418.286
$ws.Range("N41").Value = -11198.286
$ws.Range("H127").Value = 80911.42999999999
$ws.Range("J127").Value = 84331.664
$ws.Range("L127").Value = 84331.664
$ws.Range("N127").Value = -94251.664
$ws.Range("H132").Value = 1552.75
$ws.Range("I132").Value = 1532.1111
$ws.Range("J132").Value = 1614.6666
$ws.Range("K132").Value = 4596.3333
$ws.Range("L132").Value = 4843.9998
$ws.Range("M132").Value = -2066.3333
$ws.Range("N132").Value = -9903.9998
